# Adds the 2022-Q3 sheet to the workbook:
#  - the existing "2022-Q2" sheet is renamed to "2022-Q3" and gets new data
#  - a fresh copy of the (now renamed) sheet is placed right after it and
#    renamed back to "2022-Q2", preserving the original Q2 data untouched
#  - the "总计" summary sheet gets a new row for 2022-Q3 (inserted above the
#    existing 2022-Q2 summary row)

$wb = $excel.ActiveWorkbook

# Helper: force a cell to be stored as text (keeps values like "40.90" or
# "0.0314" from being coerced into numbers and losing their formatting),
# then drop back to the default "Normal" style so no stray number-format /
# quote-prefix style sticks to the cell.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Summary sheet ("总计"): shift the old 2022-Q2 row down and add the
#    new 2022-Q3 row above it.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Preserve formatting of the existing data row (A2 uses style index 2) by
# copying it down to row 3 before overwriting the values.
$summaryA2 = $summary.Cells.Item(2,1)
$summaryA3 = $summary.Cells.Item(3,1)
$summaryA2.Copy($summaryA3)
$summaryA3.Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q2"
$summary.Cells.Item(3,3).Value = 13
$summary.Cells.Item(3,4).Value = 1.88

# Overwrite row 2 with the new 2022-Q3 totals.
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 4
$summary.Cells.Item(2,4).Value = 1.17

# ---------------------------------------------------------------------
# 2. Rename the current "2022-Q2" sheet to "2022-Q3", duplicate it right
#    after itself (keeping the original Q2 numbers intact) and rename the
#    duplicate back to "2022-Q2".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3.Copy($null, $q3)
$q2 = $wb.Worksheets.Item(3)
$q2.Name = "2022-Q2"

# ---------------------------------------------------------------------
# 3. Replace the contents of the "2022-Q3" sheet with the new fund data.
# ---------------------------------------------------------------------
$q3.Cells.Clear()

# Header row (bold + bordered, matching the style used elsewhere in the
# workbook for table headers).
$h1 = $q3.Cells.Item(1,2)
$h2 = $q3.Cells.Item(1,3)
$h3 = $q3.Cells.Item(1,4)
$h4 = $q3.Cells.Item(1,5)
$h5 = $q3.Cells.Item(1,6)
$h6 = $q3.Cells.Item(1,7)
$h7 = $q3.Cells.Item(1,8)
Set-TextValue $h1 "基金代码"
Set-TextValue $h2 "基金名称"
Set-TextValue $h3 "基金规模"
Set-TextValue $h4 "股票总仓位"
Set-TextValue $h5 "仓位占比"
Set-TextValue $h6 "持有市值(亿元)"
Set-TextValue $h7 "仓位排名"

$summaryHeader = $summary.Cells.Item(1,2)
$summaryHeader.Copy()
$q3HeaderRange = $q3.Range("B1:H1")
$q3HeaderRange.PasteSpecial(-4122)

$q3Data = @(
    @("320003", "诺安先锋混合A",             "40.90", "76.11", "2.78", "1.1370", 7),
    @("012621", "诺安先锋混合C",             "1.13",  "76.11", "2.78", "0.0314", 7),
    @("003308", "中信建投睿利灵活配置混合A", "0.07",  "93.78", "4.19", "0.0029", 6),
    @("004635", "中信建投睿利灵活配置混合C", "0.03",  "93.78", "4.19", "0.0013", 6)
)

$row = 2
foreach ($fund in $q3Data) {
    $cellCode = $q3.Cells.Item($row,2)
    $cellName = $q3.Cells.Item($row,3)
    $cellSize = $q3.Cells.Item($row,4)
    $cellStockPos = $q3.Cells.Item($row,5)
    $cellPosPct = $q3.Cells.Item($row,6)
    $cellMktVal = $q3.Cells.Item($row,7)
    $cellRank = $q3.Cells.Item($row,8)

    $fundCode = $fund[0]
    $fundName = $fund[1]
    $fundSize = $fund[2]
    $fundStockPos = $fund[3]
    $fundPosPct = $fund[4]
    $fundMktVal = $fund[5]
    $fundRank = $fund[6]

    Set-TextValue $cellCode $fundCode
    Set-TextValue $cellName $fundName
    Set-TextValue $cellSize $fundSize
    Set-TextValue $cellStockPos $fundStockPos
    Set-TextValue $cellPosPct $fundPosPct
    Set-TextValue $cellMktVal $fundMktVal
    $cellRank.Value = $fundRank

    $row = $row + 1
}

# Column A (row index 0..3, style index 2) for the data rows.
$summaryA2b = $summary.Cells.Item(2,1)
$summaryA2b.Copy()
$q3ColARange = $q3.Range("A2:A5")
$q3ColARange.PasteSpecial(-4122)

$a2 = $q3.Cells.Item(2,1)
$a3 = $q3.Cells.Item(3,1)
$a4 = $q3.Cells.Item(4,1)
$a5 = $q3.Cells.Item(5,1)
$a2.Value = 0
$a3.Value = 1
$a4.Value = 2
$a5.Value = 3
